# Update the Week1 timesheet with rows 10 and 11 of logged work
# (spreadsheet rows 31-35), matching the new entries added to the
# October 2021 timesheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week1")

# --- Row 31 : SlNo 10, 10/14/2021, camunda / Camunda Basic, 2.5 hrs ---
$ws.Range("B31").Value = 10
$ws.Range("C31").Value = 44483
$ws.Range("D31").Value = "camunda"
$ws.Range("E31").Value = "Camunda Basic"
$ws.Range("F31").Value = 2.5

# --- Row 32 : java / collection, 3 hrs ---
$ws.Range("D32").Value = "java"
$ws.Range("E32").Value = "collection"
$ws.Range("F32").Value = 3

# --- Row 33 : Team Meeting / presentation session, 1.5 hrs ---
$ws.Range("D33").Value = "Team Meeting"
$ws.Range("E33").Value = "presentation session"
$ws.Range("F33").Value = 1.5

# --- Row 34 : SlNo 11, 10/15/2021, camunda / camunda setup, 2 hrs ---
$ws.Range("B34").Value = 11
$ws.Range("C34").Value = 44484
$ws.Range("D34").Value = "camunda"
$ws.Range("E34").Value = "camunda setup"
$ws.Range("F34").Value = 2

# --- Row 35 : java / java 8, 2 hrs ---
$ws.Range("D35").Value = "java"
$ws.Range("E35").Value = "java 8"
$ws.Range("F35").Value = 2

# Update the visible window / selection to match the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("F35").Select()
